# Merge the three runs of the 4th paragraph in the "Content Placeholder 2"
# shape on slide 43 into a single run with unified text, matching the
# formatting (lang="en-US" dirty="0") of the first run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(43)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$newText = "Swap design (P3239) still has a few rough edges; we are happy to bring a complete design and wording to a future telecom."

# The 4th paragraph starts right after the first three paragraphs
# (each of which is terminated by a paragraph-break character that
# counts toward TextRange character offsets).
$para1 = "We have an implementation in libc++"
$para2 = "The language feature (P2786) has been forwarded from EWG and is almost through Core."
$para3 = "High-level library relocation (P2967) is ready to be seen by LEWGI/LEWG"

$start = $para1.Length + 1 + $para2.Length + 1 + $para3.Length + 1 + 1
$len = $newText.Length

$target = $tr.Characters($start, $len)
$target.Text = $newText
